# Updated cryptos list on Mon Jul 10 03:47:24 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.177.80"
$ws.Range("E2").Value = "  -0.47%  "
$ws.Range("D3").Value = "1.863.56"
$ws.Range("E3").Value = "  -0.32%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "233.94"
$ws.Range("E5").Value = "  -0.99%  "
$ws.Range("E6").Value = "  +0.02%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4661"
$ws.Range("E7").Value = "  -1.03%  "
$ws.Range("B8").Value = "OKB"
$ws.Range("C8").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "42.85"
$ws.Range("E8").Value = "  -0.24%  "
$ws.Range("B9").Value = "Cardano"
$ws.Range("C9").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2844"
$ws.Range("E9").Value = "  -1.92%  "
$ws.Range("B10").Value = "Dogecoin"
$ws.Range("C10").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06478"
$ws.Range("E10").Value = "  -2.18%  "
$ws.Range("B11").Value = "Solana"
$ws.Range("C11").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "21.27"
$ws.Range("E11").Value = "  -2.17%  "
$ws.Range("B12").Value = "TRON"
$ws.Range("C12").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07739"
$ws.Range("E12").Value = "  -3.44%  "
$ws.Range("D13").Value = "1.881.02"
$ws.Range("E13").Value = "  +0.24%  "
$ws.Range("B14").Value = "Litecoin"
$ws.Range("C14").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "93.64"
$ws.Range("E14").Value = "  -3.79%  "
$ws.Range("B15").Value = "Polygon"
$ws.Range("C15").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6822"
$ws.Range("E15").Value = "  -1.03%  "
$ws.Range("B16").Value = "Polkadot"
$ws.Range("C16").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.037"
$ws.Range("E16").Value = "  -2.33%  "
$ws.Range("B17").Value = "BitcoinCash"
$ws.Range("C17").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "270.44"
$ws.Range("E17").Value = "  -1.36%  "
$ws.Range("B18").Value = "WrappedBTC"
$ws.Range("C18").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D18").Value = "30.154.06"
$ws.Range("E18").Value = "  -0.53%  "
$ws.Range("B19").Value = "Avalanche"
$ws.Range("C19").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.31"
$ws.Range("E19").Value = "  -5.30%  "
$ws.Range("B20").Value = "ShibaInu"
$ws.Range("C20").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007576"
$ws.Range("E20").Value = "  -1.93%  "
$ws.Range("B21").Value = "Dai"
$ws.Range("C21").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.000"
$ws.Range("E21").Value = "  +0.01%  "
$ws.Range("E22").Value = "  -0.03%  "
$ws.Range("B23").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C23").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D23").Value = "2.067.80"
$ws.Range("E23").Value = "  -2.61%  "
$ws.Range("B24").Value = "Uniswap"
$ws.Range("C24").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.131"
$ws.Range("E24").Value = "  -3.51%  "
$ws.Range("B25").Value = "Chainlink"
$ws.Range("C25").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "6.086"
$ws.Range("E25").Value = "  -2.23%  "
$ws.Range("B26").Value = "Cosmos"
$ws.Range("C26").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.338"
$ws.Range("E26").Value = "  +0.59%  "
$ws.Range("B27").Value = "Monero"
$ws.Range("C27").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "165.44"
$ws.Range("E27").Value = "  -1.32%  "
$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.55"
$ws.Range("E28").Value = "  -2.33%  "
$ws.Range("B29").Value = "LidoDAOToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.888"
$ws.Range("E29").Value = "  -3.87%  "
$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.363"
$ws.Range("E30").Value = "  -0.44%  "
$ws.Range("B31").Value = "Stellar"
$ws.Range("C31").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09832"
$ws.Range("E31").Value = "  -1.19%  "
$ws.Range("B32").Value = "PancakeSwap"
$ws.Range("C32").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.450"
$ws.Range("E32").Value = "  -0.87%  "
$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.224"
$ws.Range("E33").Value = "  -3.40%  "
$ws.Range("B34").Value = "InternetComputer(DFINITY)"
$ws.Range("C34").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.988"
$ws.Range("E34").Value = "  -2.60%  "
$ws.Range("B35").Value = "Hedera"
$ws.Range("C35").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.04665"
$ws.Range("E35").Value = "  -0.84%  "
$ws.Range("B36").Value = "ARBITRUM"
$ws.Range("C36").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.118"
$ws.Range("E36").Value = "  -1.41%  "
$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.6891"
$ws.Range("E37").Value = "  -2.07%  "
$ws.Range("B38").Value = "HuobiToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.709"
$ws.Range("E38").Value = "  +0.23%  "
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01837"
$ws.Range("E39").Value = "  -2.39%  "
$ws.Range("B40").Value = "MXToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.751"
$ws.Range("E40").Value = "  +4.34%  "
$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.302"
$ws.Range("E41").Value = "  -0.19%  "
$ws.Range("B42").Value = "Aave"
$ws.Range("C42").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "70.82"
$ws.Range("E42").Value = "  -3.55%  "
$ws.Range("B43").Value = "PaxDollar"
$ws.Range("C43").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.9999"
$ws.Range("E43").Value = "  +0.05%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.8324"
$ws.Range("E44").Value = "  -1.18%  "
$ws.Range("B45").Value = "RenderToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.880"
$ws.Range("E45").Value = "  -4.38%  "
$ws.Range("B46").Value = "Quant"
$ws.Range("C46").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "102.02"
$ws.Range("E46").Value = "  -1.51%  "
$ws.Range("B47").Value = "TheSandbox"
$ws.Range("C47").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4056"
$ws.Range("E47").Value = "  -2.78%  "
$ws.Range("B48").Value = "Maker"
$ws.Range("C48").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "934.38"
$ws.Range("E48").Value = "  +0.23%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.945"
$ws.Range("E49").Value = "  -2.76%  "
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.021"
$ws.Range("E50").Value = "  -2.35%  "
$ws.Range("B51").Value = "Elrond"
$ws.Range("C51").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "34.01"
$ws.Range("E51").Value = "  -1.45%  "
